$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.02261533333333333
$ws.Cells.Item(2, 8).Value = 0.067846
$ws.Cells.Item(2, 9).Value = 0.0136300417804998
$ws.Cells.Item(2, 10).Value = 0.0136300417804998
$ws.Cells.Item(2, 13).Value = 2.325008666666667
$ws.Cells.Item(2, 14).Value = 6.975026
$ws.Cells.Item(2, 15).Value = 0.05445297772988467
$ws.Cells.Item(2, 16).Value = 0.05445297772988466
$ws.Cells.Item(2, 17).Value = 0.05258084599955556
$ws.Cells.Item(2, 18).Value = 0.473227613996
$ws.Cells.Item(2, 19).Value = 0.000742196361530953
$ws.Cells.Item(2, 20).Value = 0.0007421963615309529

$ws.Cells.Item(3, 7).Value = 0.02261533333333333
$ws.Cells.Item(3, 8).Value = 0.067846
$ws.Cells.Item(3, 9).Value = 0.0136300417804998
$ws.Cells.Item(3, 10).Value = 0.0136300417804998
$ws.Cells.Item(3, 15).Value = 0.4529132218878514
$ws.Cells.Item(3, 16).Value = 0.4529132218878514
$ws.Cells.Item(3, 17).Value = 0.4373417462931112
$ws.Cells.Item(3, 18).Value = 3.936075716638
$ws.Cells.Item(3, 19).Value = 0.00617322613727219
$ws.Cells.Item(3, 20).Value = 0.00617322613727219

$ws.Cells.Item(4, 7).Value = 0.02261533333333333
$ws.Cells.Item(4, 8).Value = 0.067846
$ws.Cells.Item(4, 9).Value = 0.0136300417804998
$ws.Cells.Item(4, 10).Value = 0.0136300417804998
$ws.Cells.Item(4, 15).Value = 0.492633800382264
$ws.Cells.Item(4, 16).Value = 0.492633800382264
$ws.Cells.Item(4, 17).Value = 0.4756967033202222
$ws.Cells.Item(4, 18).Value = 4.281270329882
$ws.Cells.Item(4, 19).Value = 0.006714619281696655
$ws.Cells.Item(4, 20).Value = 0.006714619281696655

$ws.Cells.Item(5, 9).Value = 0.7621926756656362
$ws.Cells.Item(5, 10).Value = 0.7621926756656363
$ws.Cells.Item(5, 13).Value = 2.325008666666667
$ws.Cells.Item(5, 14).Value = 6.975026
$ws.Cells.Item(5, 15).Value = 0.05445297772988467
$ws.Cells.Item(5, 16).Value = 0.05445297772988466
$ws.Cells.Item(5, 17).Value = 2.940323760305778
$ws.Cells.Item(5, 18).Value = 26.462913842752
$ws.Cells.Item(5, 19).Value = 0.0415036607939021
$ws.Cells.Item(5, 20).Value = 0.04150366079390209

$ws.Cells.Item(6, 9).Value = 0.7621926756656362
$ws.Cells.Item(6, 10).Value = 0.7621926756656363
$ws.Cells.Item(6, 15).Value = 0.4529132218878514
$ws.Cells.Item(6, 16).Value = 0.4529132218878514
$ws.Cells.Item(6, 19).Value = 0.3452071404350455
$ws.Cells.Item(6, 20).Value = 0.3452071404350455

$ws.Cells.Item(7, 9).Value = 0.7621926756656362
$ws.Cells.Item(7, 10).Value = 0.7621926756656363
$ws.Cells.Item(7, 15).Value = 0.492633800382264
$ws.Cells.Item(7, 16).Value = 0.492633800382264
$ws.Cells.Item(7, 19).Value = 0.3754818744366887
$ws.Cells.Item(7, 20).Value = 0.3754818744366887

$ws.Cells.Item(8, 9).Value = 0.2241772825538639
$ws.Cells.Item(8, 10).Value = 0.224177282553864
$ws.Cells.Item(8, 13).Value = 2.325008666666667
$ws.Cells.Item(8, 14).Value = 6.975026
$ws.Cells.Item(8, 15).Value = 0.05445297772988467
$ws.Cells.Item(8, 16).Value = 0.05445297772988466
$ws.Cells.Item(8, 17).Value = 0.864812548662
$ws.Cells.Item(8, 18).Value = 7.783312937958
$ws.Cells.Item(8, 19).Value = 0.01220712057445162
$ws.Cells.Item(8, 20).Value = 0.01220712057445161

$ws.Cells.Item(9, 9).Value = 0.2241772825538639
$ws.Cells.Item(9, 10).Value = 0.224177282553864
$ws.Cells.Item(9, 15).Value = 0.4529132218878514
$ws.Cells.Item(9, 16).Value = 0.4529132218878514
$ws.Cells.Item(9, 19).Value = 0.1015328553155337
$ws.Cells.Item(9, 20).Value = 0.1015328553155337

$ws.Cells.Item(10, 9).Value = 0.2241772825538639
$ws.Cells.Item(10, 10).Value = 0.224177282553864
$ws.Cells.Item(10, 15).Value = 0.492633800382264
$ws.Cells.Item(10, 16).Value = 0.492633800382264
$ws.Cells.Item(10, 19).Value = 0.1104373066638786
$ws.Cells.Item(10, 20).Value = 0.1104373066638786
